$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Clear existing rows that are being fully rewritten
$ws.Range("A98:AY98").ClearContents()
$ws.Range("A99:AY99").ClearContents()
$ws.Range("A100:AY100").ClearContents()
$ws.Range("A101:AY101").ClearContents()
$ws.Range("A102:AY102").ClearContents()
$ws.Range("A103:AY103").ClearContents()
$ws.Range("A104:AY104").ClearContents()
$ws.Range("A105:AY105").ClearContents()
$ws.Range("A106:AY106").ClearContents()
$ws.Range("A109:AY109").ClearContents()
$ws.Range("A110:AY110").ClearContents()
$ws.Range("A111:AY111").ClearContents()
$ws.Range("A112:AY112").ClearContents()

# Row 98
$ws.Range("A98").Value = 111791986
$ws.Range("B98").Value = 88982
$ws.Range("C98").Value = "Ovaliderad"
$ws.Range("D98").Value = "NT"
$ws.Range("E98").Value = 937
$ws.Range("F98").Value = "Vit vedfingersvamp"
$ws.Range("G98").Value = "Lentaria epichnoa"
$ws.Range("H98").Value = "(Fr.) Corner"
$ws.Range("P98").Value = "Storvreta, Upl"
$ws.Range("Q98").Value = 650060.7258570738
$ws.Range("R98").Value = 6648731.505709799
$ws.Range("S98").Value = 25
$ws.Range("T98").Value = "Uppsala"
$ws.Range("U98").Value = "Uppsala"
$ws.Range("V98").Value = "Uppland"
$ws.Range("W98").Value = "Uppsala"
$ws.Range("Y98").NumberFormat = "@"
$ws.Range("Y98").Value = "2022-10-01"
$ws.Range("Z98").Value = "19:25"
$ws.Range("AA98").NumberFormat = "@"
$ws.Range("AA98").Value = "2022-10-01"
$ws.Range("AB98").Value = "19:25"
$ws.Range("AC98").Value = "Noterad ett år sent efter att ha problem med att rapportera! Två kraftigt rötade asplågor intill stig löpande N->S mellan Himmelsvägen och hästgården/huset efter skogen. Delar lågor med biskopsmössor."
$ws.Range("AD98").Value = $false
$ws.Range("AE98").Value = $false
$ws.Range("AG98").Value = $false
$ws.Range("AW98").Value = "Henrik Lysell"
$ws.Range("AX98").Value = "Henrik Lysell"

# Row 99
$ws.Range("A99").Value = 111790412
$ws.Range("B99").Value = 88909
$ws.Range("C99").Value = "Ovaliderad"
$ws.Range("D99").Value = "VU"
$ws.Range("E99").Value = 720
$ws.Range("F99").Value = "Violgubbe"
$ws.Range("G99").Value = "Gomphus clavatus"
$ws.Range("H99").Value = "(Pers.) Gray"
$ws.Range("P99").Value = "Storvreta, Upl"
$ws.Range("Q99").Value = 650060.7258570738
$ws.Range("R99").Value = 6648731.505709799
$ws.Range("S99").Value = 25
$ws.Range("T99").Value = "Uppsala"
$ws.Range("U99").Value = "Uppsala"
$ws.Range("V99").Value = "Uppland"
$ws.Range("W99").Value = "Uppsala"
$ws.Range("Y99").NumberFormat = "@"
$ws.Range("Y99").Value = "2023-08-30"
$ws.Range("Z99").Value = "18:04"
$ws.Range("AA99").NumberFormat = "@"
$ws.Range("AA99").Value = "2023-08-30"
$ws.Range("AB99").Value = "18:04"
$ws.Range("AC99").Value = "Separat mycel med många fruktkroppar intill ett annat stort mycel. Båda intill ett kärr."
$ws.Range("AD99").Value = $false
$ws.Range("AE99").Value = $false
$ws.Range("AG99").Value = $false
$ws.Range("AW99").Value = "Henrik Lysell"
$ws.Range("AX99").Value = "Henrik Lysell"

# Row 100
$ws.Range("A100").Value = 111790550
$ws.Range("B100").Value = 90655
$ws.Range("C100").Value = "Ovaliderad"
$ws.Range("D100").Value = "VU"
$ws.Range("E100").Value = 150
$ws.Range("F100").Value = "Grangråticka"
$ws.Range("G100").Value = "Boletopsis leucomelaena"
$ws.Range("H100").Value = "(Pers.) Fayod"
$ws.Range("P100").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q100").Value = 649935.1292544806
$ws.Range("R100").Value = 6648620.078297745
$ws.Range("S100").Value = 20
$ws.Range("T100").Value = "Uppsala"
$ws.Range("U100").Value = "Uppsala"
$ws.Range("V100").Value = "Uppland"
$ws.Range("W100").Value = "Uppsala"
$ws.Range("Y100").NumberFormat = "@"
$ws.Range("Y100").Value = "2023-08-30"
$ws.Range("Z100").Value = "18:21"
$ws.Range("AA100").NumberFormat = "@"
$ws.Range("AA100").Value = "2023-08-30"
$ws.Range("AB100").Value = "18:21"
$ws.Range("AD100").Value = $false
$ws.Range("AE100").Value = $false
$ws.Range("AG100").Value = $false
$ws.Range("AW100").Value = "Henrik Lysell"
$ws.Range("AX100").Value = "Henrik Lysell"

# Row 101
$ws.Range("A101").Value = 111789477
$ws.Range("B101").Value = 88909
$ws.Range("C101").Value = "Ovaliderad"
$ws.Range("D101").Value = "VU"
$ws.Range("E101").Value = 720
$ws.Range("F101").Value = "Violgubbe"
$ws.Range("G101").Value = "Gomphus clavatus"
$ws.Range("H101").Value = "(Pers.) Gray"
$ws.Range("P101").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q101").Value = 649983.3633731108
$ws.Range("R101").Value = 6648731.373486301
$ws.Range("S101").Value = 10
$ws.Range("T101").Value = "Uppsala"
$ws.Range("U101").Value = "Uppsala"
$ws.Range("V101").Value = "Uppland"
$ws.Range("W101").Value = "Uppsala"
$ws.Range("Y101").NumberFormat = "@"
$ws.Range("Y101").Value = "2023-08-30"
$ws.Range("Z101").Value = "17:09"
$ws.Range("AA101").NumberFormat = "@"
$ws.Range("AA101").Value = "2023-08-30"
$ws.Range("AB101").Value = "17:09"
$ws.Range("AC101").Value = "Mellan stenblock ca 50 m bort ifrån föregående rapporterade mycel. Bland ett myller av granlågor."
$ws.Range("AD101").Value = $false
$ws.Range("AE101").Value = $false
$ws.Range("AG101").Value = $false
$ws.Range("AW101").Value = "Henrik Lysell"
$ws.Range("AX101").Value = "Henrik Lysell"

# Row 102
$ws.Range("A102").Value = 111789319
$ws.Range("B102").Value = 88915
$ws.Range("C102").Value = "Ovaliderad"
$ws.Range("D102").Value = "NT"
$ws.Range("E102").Value = 5734
$ws.Range("F102").Value = "Druvfingersvamp"
$ws.Range("G102").Value = "Ramaria botrytis"
$ws.Range("H102").Value = "(Pers.:Fr.) Bourdot"
$ws.Range("P102").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q102").Value = 650042.1747608959
$ws.Range("R102").Value = 6648755.327700124
$ws.Range("S102").Value = 10
$ws.Range("T102").Value = "Uppsala"
$ws.Range("U102").Value = "Uppsala"
$ws.Range("V102").Value = "Uppland"
$ws.Range("W102").Value = "Uppsala"
$ws.Range("Y102").NumberFormat = "@"
$ws.Range("Y102").Value = "2023-08-30"
$ws.Range("Z102").Value = "16:58"
$ws.Range("AA102").NumberFormat = "@"
$ws.Range("AA102").Value = "2023-08-30"
$ws.Range("AB102").Value = "16:58"
$ws.Range("AC102").Value = "Mitt i en stigkorsning"
$ws.Range("AD102").Value = $false
$ws.Range("AE102").Value = $false
$ws.Range("AG102").Value = $false
$ws.Range("AH102").Value = "Blåbärsgranskog"
$ws.Range("AW102").Value = "Henrik Lysell"
$ws.Range("AX102").Value = "Henrik Lysell"

# Row 103
$ws.Range("A103").Value = 111790914
$ws.Range("B103").Value = 90662
$ws.Range("C103").Value = "Ovaliderad"
$ws.Range("D103").Value = "LC"
$ws.Range("E103").Value = 4363
$ws.Range("F103").Value = "Zontaggsvamp"
$ws.Range("G103").Value = "Hydnellum concrescens"
$ws.Range("H103").Value = "(Pers.) Banker"
$ws.Range("P103").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q103").Value = 650131.1687508342
$ws.Range("R103").Value = 6648444.465040453
$ws.Range("S103").Value = 100
$ws.Range("T103").Value = "Uppsala"
$ws.Range("U103").Value = "Uppsala"
$ws.Range("V103").Value = "Uppland"
$ws.Range("W103").Value = "Uppsala"
$ws.Range("Y103").NumberFormat = "@"
$ws.Range("Y103").Value = "2023-08-30"
$ws.Range("Z103").Value = "18:43"
$ws.Range("AA103").NumberFormat = "@"
$ws.Range("AA103").Value = "2023-08-30"
$ws.Range("AB103").Value = "18:43"
$ws.Range("AD103").Value = $false
$ws.Range("AE103").Value = $false
$ws.Range("AG103").Value = $false
$ws.Range("AW103").Value = "Henrik Lysell"
$ws.Range("AX103").Value = "Henrik Lysell"

# Row 104
$ws.Range("A104").Value = 111789261
$ws.Range("B104").Value = 88915
$ws.Range("C104").Value = "Ovaliderad"
$ws.Range("D104").Value = "NT"
$ws.Range("E104").Value = 5734
$ws.Range("F104").Value = "Druvfingersvamp"
$ws.Range("G104").Value = "Ramaria botrytis"
$ws.Range("H104").Value = "(Pers.:Fr.) Bourdot"
$ws.Range("P104").Value = "Storvreta, Upl"
$ws.Range("Q104").Value = 650060.7258570738
$ws.Range("R104").Value = 6648731.505709799
$ws.Range("S104").Value = 25
$ws.Range("T104").Value = "Uppsala"
$ws.Range("U104").Value = "Uppsala"
$ws.Range("V104").Value = "Uppland"
$ws.Range("W104").Value = "Uppsala"
$ws.Range("Y104").NumberFormat = "@"
$ws.Range("Y104").Value = "2023-08-30"
$ws.Range("Z104").Value = "17:02"
$ws.Range("AA104").NumberFormat = "@"
$ws.Range("AA104").Value = "2023-08-30"
$ws.Range("AB104").Value = "17:02"
$ws.Range("AD104").Value = $false
$ws.Range("AE104").Value = $false
$ws.Range("AG104").Value = $false
$ws.Range("AW104").Value = "Henrik Lysell"
$ws.Range("AX104").Value = "Henrik Lysell"

# Row 105
$ws.Range("A105").Value = 111789785
$ws.Range("B105").Value = 90187
$ws.Range("C105").Value = "Ovaliderad"
$ws.Range("D105").Value = "NT"
$ws.Range("E105").Value = 2014
$ws.Range("F105").Value = "Koralltaggsvamp"
$ws.Range("G105").Value = "Hericium coralloides"
$ws.Range("H105").Value = "(Scop.:Fr.) Pers."
$ws.Range("P105").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q105").Value = 649874.9776032839
$ws.Range("R105").Value = 6648703.402536019
$ws.Range("S105").Value = 20
$ws.Range("T105").Value = "Uppsala"
$ws.Range("U105").Value = "Uppsala"
$ws.Range("V105").Value = "Uppland"
$ws.Range("W105").Value = "Uppsala"
$ws.Range("Y105").NumberFormat = "@"
$ws.Range("Y105").Value = "2023-08-30"
$ws.Range("Z105").Value = "17:35"
$ws.Range("AA105").NumberFormat = "@"
$ws.Range("AA105").Value = "2023-08-30"
$ws.Range("AB105").Value = "17:35"
$ws.Range("AC105").Value = "I starkt rötad asplåga. 2 fruktkroppar. Ca 70 meter ifrån skogsbrynet till en elledningsgata. Svårt att få exakt position pga mobiltäckning."
$ws.Range("AD105").Value = $false
$ws.Range("AE105").Value = $false
$ws.Range("AG105").Value = $false
$ws.Range("AW105").Value = "Henrik Lysell"
$ws.Range("AX105").Value = "Henrik Lysell"

# Row 106
$ws.Range("A106").Value = 111789368
$ws.Range("B106").Value = 88909
$ws.Range("C106").Value = "Ovaliderad"
$ws.Range("D106").Value = "VU"
$ws.Range("E106").Value = 720
$ws.Range("F106").Value = "Violgubbe"
$ws.Range("G106").Value = "Gomphus clavatus"
$ws.Range("H106").Value = "(Pers.) Gray"
$ws.Range("P106").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q106").Value = 650001.3452252811
$ws.Range("R106").Value = 6648758.685092625
$ws.Range("S106").Value = 10
$ws.Range("T106").Value = "Uppsala"
$ws.Range("U106").Value = "Uppsala"
$ws.Range("V106").Value = "Uppland"
$ws.Range("W106").Value = "Uppsala"
$ws.Range("Y106").NumberFormat = "@"
$ws.Range("Y106").Value = "2023-08-30"
$ws.Range("Z106").Value = "17:09"
$ws.Range("AA106").NumberFormat = "@"
$ws.Range("AA106").Value = "2023-08-30"
$ws.Range("AB106").Value = "17:09"
$ws.Range("AC106").Value = "I ett ca 5 meter långt stråk med många fruktkroppar."
$ws.Range("AD106").Value = $false
$ws.Range("AE106").Value = $false
$ws.Range("AG106").Value = $false
$ws.Range("AH106").Value = "Blåbärsgranskog"
$ws.Range("AW106").Value = "Henrik Lysell"
$ws.Range("AX106").Value = "Henrik Lysell"

# Row 109
$ws.Range("A109").Value = 111984708
$ws.Range("B109").Value = 88909
$ws.Range("C109").Value = "Ovaliderad"
$ws.Range("D109").Value = "VU"
$ws.Range("E109").Value = 720
$ws.Range("F109").Value = "Violgubbe"
$ws.Range("G109").Value = "Gomphus clavatus"
$ws.Range("H109").Value = "(Pers.) Gray"
$ws.Range("P109").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q109").Value = 650056.8292729721
$ws.Range("R109").Value = 6648629.533134428
$ws.Range("S109").Value = 20
$ws.Range("T109").Value = "Uppsala"
$ws.Range("U109").Value = "Uppsala"
$ws.Range("V109").Value = "Uppland"
$ws.Range("W109").Value = "Uppsala"
$ws.Range("Y109").NumberFormat = "@"
$ws.Range("Y109").Value = "2023-09-09"
$ws.Range("Z109").Value = "14:54"
$ws.Range("AA109").NumberFormat = "@"
$ws.Range("AA109").Value = "2023-09-09"
$ws.Range("AB109").Value = "14:54"
$ws.Range("AC109").Value = "Under granar, precis intill kärret och en halvmeter från stig. Två samlingar."
$ws.Range("AD109").Value = $false
$ws.Range("AE109").Value = $false
$ws.Range("AG109").Value = $false
$ws.Range("AW109").Value = "Henrik Lysell"
$ws.Range("AX109").Value = "Henrik Lysell"

# Row 110
$ws.Range("A110").Value = 111984394
$ws.Range("B110").Value = 90662
$ws.Range("C110").Value = "Ovaliderad"
$ws.Range("D110").Value = "LC"
$ws.Range("E110").Value = 4363
$ws.Range("F110").Value = "Zontaggsvamp"
$ws.Range("G110").Value = "Hydnellum concrescens"
$ws.Range("H110").Value = "(Pers.) Banker"
$ws.Range("P110").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q110").Value = 650056.8292729721
$ws.Range("R110").Value = 6648629.533134428
$ws.Range("S110").Value = 20
$ws.Range("T110").Value = "Uppsala"
$ws.Range("U110").Value = "Uppsala"
$ws.Range("V110").Value = "Uppland"
$ws.Range("W110").Value = "Uppsala"
$ws.Range("Y110").NumberFormat = "@"
$ws.Range("Y110").Value = "2023-09-09"
$ws.Range("Z110").Value = "13:16"
$ws.Range("AA110").NumberFormat = "@"
$ws.Range("AA110").Value = "2023-09-09"
$ws.Range("AB110").Value = "13:16"
$ws.Range("AC110").Value = "Nedanför granbacke i en stig (kärr ca 40 meter NV)"
$ws.Range("AD110").Value = $false
$ws.Range("AE110").Value = $false
$ws.Range("AG110").Value = $false
$ws.Range("AW110").Value = "Henrik Lysell"
$ws.Range("AX110").Value = "Henrik Lysell"

# Row 111
$ws.Range("A111").Value = 112063812
$ws.Range("B111").Value = 90678
$ws.Range("C111").Value = "Ovaliderad"
$ws.Range("D111").Value = "LC"
$ws.Range("E111").Value = 4366
$ws.Range("F111").Value = "Skarp dropptaggsvamp"
$ws.Range("G111").Value = "Hydnellum peckii"
$ws.Range("H111").Value = "Banker"
$ws.Range("I111").NumberFormat = "@"
$ws.Range("I111").Value = "4"
$ws.Range("P111").Value = "Barrskog V Himmelsvägen, Upl"
$ws.Range("Q111").Value = 650048.6548849966
$ws.Range("R111").Value = 6648645.250869188
$ws.Range("S111").Value = 10
$ws.Range("T111").Value = "Uppsala"
$ws.Range("U111").Value = "Uppsala"
$ws.Range("V111").Value = "Uppland"
$ws.Range("W111").Value = "Uppsala"
$ws.Range("Y111").NumberFormat = "@"
$ws.Range("Y111").Value = "2023-09-13"
$ws.Range("Z111").Value = "00:00"
$ws.Range("AA111").NumberFormat = "@"
$ws.Range("AA111").Value = "2023-09-13"
$ws.Range("AB111").Value = "00:00"
$ws.Range("AC111").Value = "Flera små, ganska torra, fruktkroppar"
$ws.Range("AD111").Value = $false
$ws.Range("AE111").Value = $false
$ws.Range("AG111").Value = $false
$ws.Range("AW111").Value = "Gry Benediktson"
$ws.Range("AX111").Value = "Gry Benediktson"

# Row 112
$ws.Range("A112").Value = 112139297
$ws.Range("B112").Value = 88909
$ws.Range("C112").Value = "Ovaliderad"
$ws.Range("D112").Value = "VU"
$ws.Range("E112").Value = 720
$ws.Range("F112").Value = "Violgubbe"
$ws.Range("G112").Value = "Gomphus clavatus"
$ws.Range("H112").Value = "(Pers.) Gray"
$ws.Range("P112").Value = "barrskog V himmelsvägen, Upl"
$ws.Range("Q112").Value = 650012.2168472378
$ws.Range("R112").Value = 6648763.138812571
$ws.Range("S112").Value = 10
$ws.Range("T112").Value = "Uppsala"
$ws.Range("U112").Value = "Uppsala"
$ws.Range("V112").Value = "Uppland"
$ws.Range("W112").Value = "Uppsala"
$ws.Range("Y112").NumberFormat = "@"
$ws.Range("Y112").Value = "2023-09-14"
$ws.Range("Z112").Value = "00:00"
$ws.Range("AA112").NumberFormat = "@"
$ws.Range("AA112").Value = "2023-09-14"
$ws.Range("AB112").Value = "00:00"
$ws.Range("AD112").Value = $false
$ws.Range("AE112").Value = $false
$ws.Range("AG112").Value = $false
$ws.Range("AW112").Value = "Gry Benediktson"
$ws.Range("AX112").Value = "Gry Benediktson"

# Row 113
$ws.Range("A113").Value = 112139372
$ws.Range("B113").Value = 90678
$ws.Range("C113").Value = "Ovaliderad"
$ws.Range("D113").Value = "LC"
$ws.Range("E113").Value = 4366
$ws.Range("F113").Value = "Skarp dropptaggsvamp"
$ws.Range("G113").Value = "Hydnellum peckii"
$ws.Range("H113").Value = "Banker"
$ws.Range("P113").Value = "barrskog V himmelsvägen, Upl"
$ws.Range("Q113").Value = 649885.7169699427
$ws.Range("R113").Value = 6648971.151017974
$ws.Range("S113").Value = 10
$ws.Range("T113").Value = "Uppsala"
$ws.Range("U113").Value = "Uppsala"
$ws.Range("V113").Value = "Uppland"
$ws.Range("W113").Value = "Uppsala"
$ws.Range("Y113").NumberFormat = "@"
$ws.Range("Y113").Value = "2023-09-14"
$ws.Range("Z113").Value = "00:00"
$ws.Range("AA113").NumberFormat = "@"
$ws.Range("AA113").Value = "2023-09-14"
$ws.Range("AB113").Value = "00:00"
$ws.Range("AC113").Value = "På/intill stigen"
$ws.Range("AD113").Value = $false
$ws.Range("AE113").Value = $false
$ws.Range("AG113").Value = $false
$ws.Range("AW113").Value = "Gry Benediktson"
$ws.Range("AX113").Value = "Gry Benediktson"
